$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 323, shifting existing rows 323:370 down to 324:371.
# Excel's Insert() propagates the row's cell formatting (e.g. the date style on
# column D) from the surrounding rows automatically.
$ws.Rows(323).Insert()

# Populate the newly inserted row 323 with this week's new price entry
# (a duplicate of the most recent "Asterix" / "1a (cosecha)" record, dated one
# week later than the previous latest entry).
$row = 323
$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item($row, 3).Value = 'Bíobío'
$ws.Cells.Item($row, 4).Value = 44918
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = 100114001
$ws.Cells.Item($row, 7).Value = 'Papa'
$ws.Cells.Item($row, 8).Value = 'Asterix'
$ws.Cells.Item($row, 9).Value = '1a (cosecha)'
$ws.Cells.Item($row, 10).Value = 10000
$ws.Cells.Item($row, 11).Value = 11000
$ws.Cells.Item($row, 12).Value = 12000
$ws.Cells.Item($row, 13).Value = 11500
$ws.Cells.Item($row, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item($row, 15).Value = 'Provincia de Arauco'
$ws.Cells.Item($row, 16).Value = 460
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = 'Hortaliza'
